$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-0.02314597604078636;  C=0.3579920056255013;  D=0.1782699060034266;  E=0.4222202103209018;  F=0.4363822494547141;  G=15},
    @{Row=3;  B=0.01177902937386835;   C=0.3233677337305382;  D=0.1387525778210147;  E=0.3724950708680783;  F=0.3863630867095333;  G=14},
    @{Row=4;  B=-0.008234150994489252; C=0.3656812685479343;  D=0.1736781030421116;  E=0.4167470492302394;  F=0.4336794057244024;  G=13},
    @{Row=5;  B=0.04114622104959433;   C=0.2438785208639479;  D=0.09301430872978532; E=0.3049824728239072;  F=0.315631483433317;   G=12},
    @{Row=6;  B=-0.004251297071472607; C=0.3196557001756482;  D=0.1664324393548997;  E=0.4079613209054257;  F=0.4278502102499434;  G=11},
    @{Row=7;  B=-0.03152289704198964;  C=0.3137577720403564;  D=0.1200567451426497;  E=0.3464920563918454;  F=0.3637200573468097;  G=10},
    @{Row=8;  B=-0.05669611638992965;  C=0.309241221828214;   D=0.12362493933599;    E=0.3516033835673229;  F=0.3680513563860018;  G=9},
    @{Row=9;  B=-0.04410238231315148;  C=0.2487243691024244;  D=0.1063280554679864;  E=0.3260798299005727;  F=0.3453909343541705;  G=8},
    @{Row=10; B=-0.07289500678705241;  C=0.3460499651212735;  D=0.2127274776992719;  E=0.4612238910759848;  F=0.4919174337551497;  G=7},
    @{Row=11; B=-0.05253816173664939;  C=0.2093953811075043;  D=0.05818608614024737; E=0.2412179225104291;  F=0.2578972532655152;  G=6}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
